# Add a new worksheet "invalidZip" after the existing "doSearchAtmLocation"
# sheet, populate it with a header + three sample zip values, select A4,
# and make it the active sheet (matching the author's "switch languages"
# test-data commit).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "invalidZip"

$newSheet.Range("A1").Value = "invalidZip"
$newSheet.Range("A2").Value = 123
$newSheet.Range("A3").Value = 456
$newSheet.Range("A4").Value = 789

[void]$newSheet.Range("A4").Select()
$newSheet.Activate()
